$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 482.8
$ws.Range("I4").Value = 482.8
$ws.Range("K4").Value = 482.8
$ws.Range("M4").Value = -368.8
$ws.Range("H18").Value = 1450
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H62").Value = 16755.875
$ws.Range("I62").Value = 23809.6
$ws.Range("K62").Value = 23809.6
$ws.Range("M62").Value = -23185.6
$ws.Range("H65").Value = 16755.875
$ws.Range("I65").Value = 23809.6
$ws.Range("K65").Value = 119048
$ws.Range("M65").Value = -115928
$ws.Range("H86").Value = 4811.1333
$ws.Range("I86").Value = 4103.3335
$ws.Range("J86").Value = 5283
$ws.Range("K86").Value = 4103.3335
$ws.Range("L86").Value = 5283
$ws.Range("M86").Value = -2980.3335
$ws.Range("N86").Value = -7529
$ws.Range("H87").Value = 27020
$ws.Range("J87").Value = 27750
$ws.Range("L87").Value = 27750
$ws.Range("N87").Value = -30246
$ws.Range("H89").Value = 4811.1333
$ws.Range("I89").Value = 4103.3335
$ws.Range("J89").Value = 5283
$ws.Range("K89").Value = 20516.6675
$ws.Range("L89").Value = 26415
$ws.Range("M89").Value = -14900.6675
$ws.Range("N89").Value = -37647
$ws.Range("H90").Value = 27020
$ws.Range("J90").Value = 27750
$ws.Range("L90").Value = 83250
$ws.Range("N90").Value = -95730
$ws.Range("H113").Value = 9473.091
$ws.Range("I113").Value = 6841
$ws.Range("K113").Value = 6841
$ws.Range("M113").Value = -3587
$ws.Range("H116").Value = 9280.700000000001
$ws.Range("I116").Value = 9858.223
$ws.Range("J116").Value = 8808.182000000001
$ws.Range("K116").Value = 9858.223
$ws.Range("L116").Value = 8808.182000000001
$ws.Range("M116").Value = -6416.223
$ws.Range("N116").Value = -15692.182
$ws.Range("H132").Value = 11855.143
$ws.Range("I132").Value = 1592.5869
$ws.Range("K132").Value = 4777.7607
$ws.Range("M132").Value = -2247.7607
$ws.Range("H138").Value = 3881.6667
$ws.Range("J138").Value = 4950
$ws.Range("L138").Value = 14850
$ws.Range("N138").Value = -25130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20413858
$ws.Range("I32").Value = 22228298
$ws.Range("J32").Value = 1415.75
$ws.Range("K32").Value = 22228298
$ws.Range("L32").Value = 1415.75
$ws.Range("M32").Value = -22228011
$ws.Range("N32").Value = -1989.75
$ws.Range("H39").Value = 4732.4614
$ws.Range("I39").Value = 3074.5715
$ws.Range("K39").Value = 3074.5715
$ws.Range("M39").Value = -2554.5715
$ws.Range("H45").Value = 4872.0386
$ws.Range("I45").Value = 3932.2222
$ws.Range("J45").Value = 6986.625
$ws.Range("K45").Value = 3932.2222
$ws.Range("L45").Value = 6986.625
$ws.Range("M45").Value = -3555.2222
$ws.Range("N45").Value = -7740.625
$ws.Range("H61").Value = 2074.4211
$ws.Range("I61").Value = 2074.4211
$ws.Range("K61").Value = 2074.4211
$ws.Range("M61").Value = -1862.4211
$ws.Range("H63").Value = 7835
$ws.Range("I63").Value = 10749.5
$ws.Range("J63").Value = 2006
$ws.Range("K63").Value = 10749.5
$ws.Range("L63").Value = 2006
$ws.Range("M63").Value = -10063.5
$ws.Range("N63").Value = -3378
$ws.Range("H66").Value = 7835
$ws.Range("I66").Value = 10749.5
$ws.Range("J66").Value = 2006
$ws.Range("K66").Value = 53747.5
$ws.Range("L66").Value = 10030
$ws.Range("M66").Value = -50315.5
$ws.Range("N66").Value = -16894
$ws.Range("H97").Value = 2964.85
$ws.Range("I97").Value = 2127.611
$ws.Range("J97").Value = 10500
$ws.Range("K97").Value = 2127.611
$ws.Range("L97").Value = 10500
$ws.Range("M97").Value = -1631.611
$ws.Range("N97").Value = -11492
$ws.Range("H136").Value = 2074.4211
$ws.Range("I136").Value = 2074.4211
$ws.Range("K136").Value = 6223.263300000001
$ws.Range("M136").Value = -3673.263300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2863
$ws.Range("I86").Value = 2295
$ws.Range("K86").Value = 2295
$ws.Range("M86").Value = -1172
$ws.Range("H89").Value = 2863
$ws.Range("I89").Value = 2295
$ws.Range("K89").Value = 11475
$ws.Range("M89").Value = -5859

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 895
$ws.Range("I21").Value = 895
$ws.Range("K21").Value = 895
$ws.Range("M21").Value = -660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 17194.928
$ws.Range("J3").Value = 19999.857
$ws.Range("L3").Value = 59999.571
$ws.Range("N3").Value = -60223.571
$ws.Range("H14").Value = 159.375
$ws.Range("I14").Value = 159.375
$ws.Range("K14").Value = 478.125
$ws.Range("M14").Value = -305.125
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H107").Value = 1499.16
$ws.Range("J107").Value = 1755.75
$ws.Range("L107").Value = 5267.25
$ws.Range("N107").Value = -9107.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 145465.38
$ws.Range("I70").Value = 228585.2
$ws.Range("J70").Value = 6932.3335
$ws.Range("K70").Value = 228585.2
$ws.Range("L70").Value = 6932.3335
$ws.Range("M70").Value = -228315.2
$ws.Range("N70").Value = -7472.3335
$ws.Range("H73").Value = 145465.38
$ws.Range("I73").Value = 228585.2
$ws.Range("J73").Value = 6932.3335
$ws.Range("K73").Value = 228585.2
$ws.Range("L73").Value = 6932.3335
$ws.Range("M73").Value = -227649.2
$ws.Range("N73").Value = -8804.333500000001
$ws.Range("H122").Value = 5686.636
$ws.Range("I122").Value = 5067.5
$ws.Range("K122").Value = 15202.5
$ws.Range("M122").Value = -12752.5
$ws.Range("H132").Value = 11118.235
$ws.Range("I132").Value = 10308.538
$ws.Range("K132").Value = 30925.614
$ws.Range("M132").Value = -28395.614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6272.727
$ws.Range("I7").Value = 6272.727
$ws.Range("K7").Value = 6272.727
$ws.Range("M7").Value = -6160.727
$ws.Range("H30").Value = 7634.8335
$ws.Range("I30").Value = 7634.8335
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 7634.8335
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -7526.8335
$ws.Range("N30").ClearContents()
$ws.Range("H40").Value = 4276.5
$ws.Range("I40").Value = 3553
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3553
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3417
$ws.Range("N40").Value = -5272
$ws.Range("H68").Value = 3305.5
$ws.Range("I68").Value = 3234.8333
$ws.Range("J68").Value = 3411.5
$ws.Range("K68").Value = 3234.8333
$ws.Range("L68").Value = 3411.5
$ws.Range("M68").Value = -2485.8333
$ws.Range("N68").Value = -4909.5
$ws.Range("H71").Value = 3305.5
$ws.Range("I71").Value = 3234.8333
$ws.Range("J71").Value = 3411.5
$ws.Range("K71").Value = 16174.1665
$ws.Range("L71").Value = 17057.5
$ws.Range("M71").Value = -12430.1665
$ws.Range("N71").Value = -24545.5
$ws.Range("H100").Value = 64168.145
$ws.Range("I100").Value = 82248.734
$ws.Range("K100").Value = 82248.734
$ws.Range("M100").Value = -81707.734
$ws.Range("H122").Value = 6373.7144
$ws.Range("I122").Value = 4922.4
$ws.Range("K122").Value = 14767.2
$ws.Range("M122").Value = -12317.2
$ws.Range("H126").Value = 6272.727
$ws.Range("I126").Value = 6272.727
$ws.Range("K126").Value = 18818.181
$ws.Range("M126").Value = -16348.181
$ws.Range("H127").Value = 45850.312
$ws.Range("J127").Value = 45850.312
$ws.Range("L127").Value = 45850.312
$ws.Range("N127").Value = -55770.312
$ws.Range("H132").Value = 4964.9536
$ws.Range("I132").Value = 5069.7437
$ws.Range("K132").Value = 15209.2311
$ws.Range("M132").Value = -12679.2311
$ws.Range("H136").Value = 2316.5833
$ws.Range("I136").Value = 2066.3333
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 6198.999899999999
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -3648.999899999999
$ws.Range("N136").Value = -12300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H99").Value = 70000
$ws.Range("J99").Value = 70000
$ws.Range("L99").Value = 70000
$ws.Range("N99").Value = -75990
$ws.Range("H100").Value = 764.125
$ws.Range("I100").Value = 780.4286
$ws.Range("K100").Value = 1560.8572
$ws.Range("M100").Value = -1019.8572
$ws.Range("H123").Value = 29785.715
$ws.Range("J123").Value = 29785.715
$ws.Range("L123").Value = 29785.715
$ws.Range("N123").Value = -39585.715
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H132").Value = 5905.875
$ws.Range("I132").Value = 5364.9
$ws.Range("J132").Value = 6807.5
$ws.Range("K132").Value = 16094.7
$ws.Range("L132").Value = 20422.5
$ws.Range("M132").Value = -13564.7
$ws.Range("N132").Value = -25482.5
$ws.Range("H135").Value = 59461.152
$ws.Range("J135").Value = 59461.152
$ws.Range("L135").Value = 59461.152
$ws.Range("N135").Value = -69601.152
$ws.Range("H138").Value = 76499.8
$ws.Range("J138").Value = 76499.8
$ws.Range("L138").Value = 76499.8
$ws.Range("N138").Value = -86779.8
